# Updated symbol list on Sat Dec 31 05:47:10 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells hold numeric-looking values stored as literal
# inline-string text in the workbook (t="inlineStr"), not real numbers.
# A leading apostrophe forces Excel to keep the exact text (so "25.50" does
# not get normalised to 25.5, "0.1340" does not become 0.134, tiny values
# like 0.00009702 do not get turned into scientific notation, etc.). Resetting
# the Style afterwards drops the implicit "Text" number-format Excel applies
# so the cell keeps its original (default) styling.

$ws.Range("D3").Value = "'25.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.097"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05571"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.474"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.016"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8187"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8440"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1340"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02854"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09380"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001511"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005962"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14OneONE"
$ws.Range("D16").Value = "'0.006116"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'3.500"
$ws.Range("D17").Style = "Normal"
$ws.Range("B20").Value = "MandalaExchangeToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D20").Value = "'0.06957"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19MandalaExchangeTokenMDX"
$ws.Range("D22").Value = "'3.766"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04718"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.001250"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'0.00009702"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002631"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003375"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.008305"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005299"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Value = "'0.002121"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
